$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on D2:E51 so numeric-looking / percent-looking strings
# are stored as text (matching the original inlineStr cell type) instead of
# being auto-converted to numbers by Excel value parsing.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '41.838.02'
$ws.Range('E2').Value = '  -1.72%  '
$ws.Range('D3').Value = '2.279.06'
$ws.Range('E3').Value = '  -2.84%  '
$ws.Range('E4').Value = '  +0.23%  '
$ws.Range('D5').Value = '316.48'
$ws.Range('E5').Value = '  +0.29%  '
$ws.Range('D6').Value = '102.12'
$ws.Range('E6').Value = '  -5.65%  '
$ws.Range('D7').Value = '0.625'
$ws.Range('E7').Value = '  -1.05%  '
$ws.Range('E8').Value = '  +0.17%  '
$ws.Range('D9').Value = '0.602'
$ws.Range('E9').Value = '  -2.22%  '
$ws.Range('D10').Value = '38.81'
$ws.Range('E10').Value = '  -6.03%  '
$ws.Range('D11').Value = '0.0904'
$ws.Range('E11').Value = '  -2.44%  '
$ws.Range('D12').Value = '8.24'
$ws.Range('E12').Value = '  -4.55%  '
$ws.Range('E13').Value = '  -0.25%  '
$ws.Range('D14').Value = '0.957'
$ws.Range('E14').Value = '  -4.02%  '
$ws.Range('D15').Value = '15.20'
$ws.Range('E15').Value = '  -4.82%  '
$ws.Range('D16').Value = '2.623.77'
$ws.Range('E16').Value = '  -2.81%  '
$ws.Range('D17').Value = '2.277.19'
$ws.Range('E17').Value = '  -2.30%  '
$ws.Range('D18').Value = '41.833.80'
$ws.Range('E18').Value = '  -1.64%  '
$ws.Range('E19').Value = '  -2.09%  '
$ws.Range('E20').Value = '  -1.25%  '
$ws.Range('D23').Value = '3.57'
$ws.Range('E23').Value = '  -1.61%  '
$ws.Range('D24').Value = '2.26'
$ws.Range('E24').Value = '  -2.35%  '
$ws.Range('D25').Value = '9.89'
$ws.Range('E25').Value = '  +5.09%  '
$ws.Range('E26').Value = '  +0.89%  '
$ws.Range('D27').Value = '10.73'
$ws.Range('E27').Value = '  -5.78%  '
$ws.Range('D28').Value = '2.30'
$ws.Range('E28').Value = '  +3.92%  '
$ws.Range('D29').Value = '22.98'
$ws.Range('E29').Value = '  +1.31%  '
$ws.Range('D30').Value = '162.45'
$ws.Range('E30').Value = '  -5.95%  '
$ws.Range('D31').Value = '0.0871'
$ws.Range('E31').Value = '  -2.31%  '
$ws.Range('D32').Value = '34.27'
$ws.Range('E32').Value = '  -7.01%  '
$ws.Range('D33').Value = '2.90'
$ws.Range('E33').Value = '  +1.18%  '
$ws.Range('D34').Value = '5.81'
$ws.Range('E34').Value = '  -3.69%  '
$ws.Range('D35').Value = '0.132'
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('E36').Value = '  -7.11%  '
$ws.Range('D37').Value = '4.57'
$ws.Range('E37').Value = '  -1.40%  '
$ws.Range('D38').Value = '2.89'
$ws.Range('E38').Value = '  +8.84%  '
$ws.Range('D39').Value = '0.0346'
$ws.Range('E39').Value = '  -4.22%  '
$ws.Range('D40').Value = '3.60'
$ws.Range('E40').Value = '  -8.40%  '
$ws.Range('D41').Value = '102.76'
$ws.Range('E41').Value = '  +20.28%  '
$ws.Range('D42').Value = '1.45'
$ws.Range('E42').Value = '  -1.28%  '
$ws.Range('D43').Value = '69.25'
$ws.Range('E43').Value = '  -1.57%  '
$ws.Range('E44').Value = '  +0.09%  '
$ws.Range('D45').Value = '0.224'
$ws.Range('E45').Value = '  -6.72%  '
$ws.Range('D46').Value = '115.17'
$ws.Range('E46').Value = '  +3.75%  '
$ws.Range('D47').Value = '11.84'
$ws.Range('E47').Value = '  -1.98%  '
$ws.Range('D48').Value = '9.00'
$ws.Range('E48').Value = '  -2.00%  '
$ws.Range('D51').Value = '1.26'
$ws.Range('E51').Value = '  -2.32%  '

# Row swaps: BitcoinCash/Litecoin (rows 21-22) and THORChain/ordi (rows 49-50)
$ws.Range('B21').Value = 'Litecoin'
$ws.Range('C21').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D21').Value = '73.41'
$ws.Range('E21').Value = '  -3.00%  '
$ws.Range('B22').Value = 'BitcoinCash'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D22').Value = '282.10'
$ws.Range('E22').Value = '  +10.48%  '
$ws.Range('B49').Value = 'ordi'
$ws.Range('C49').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D49').Value = '75.81'
$ws.Range('E49').Value = '  +1.18%  '
$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D50').Value = '5.26'
$ws.Range('E50').Value = '  -3.56%  '

# Clear the temporary NumberFormat override so the style index matches the
# original (unstyled) cells again.
$ws.Range("D2:E51").ClearFormats()
